$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("B2").Value = -11.61609041642555
$ws.Range("F2").Value = -0.3808635473251343
$ws.Range("J2").Value = -51.60860252380371
$ws.Range("K2").Value = -52.35033583641052
$ws.Range("B3").Value = -13.48009934019865
$ws.Range("F3").Value = -0.1282483339309692
$ws.Range("J3").Value = -50.4172728061676
$ws.Range("K3").Value = -50.74892771244049
$ws.Range("B4").Value = -15.17550054018648
$ws.Range("F4").Value = 0.2709441184997559
$ws.Range("J4").Value = -49.61408925056458
$ws.Range("K4").Value = -49.4236820936203
$ws.Range("B5").Value = -15.9793464043189
$ws.Range("F5").Value = 0.7460058927536011
$ws.Range("J5").Value = -49.13628923892975
$ws.Range("K5").Value = -48.54786038398743
$ws.Range("B6").Value = -16.26560880342163
$ws.Range("F6").Value = 1.053337335586548
$ws.Range("J6").Value = -48.89055645465851
$ws.Range("K6").Value = -48.0921596288681
$ws.Range("B7").Value = -16.26560880342163
$ws.Range("F7").Value = 1.053337335586548
$ws.Range("J7").Value = -48.89055645465851
$ws.Range("K7").Value = -48.0921596288681
$ws.Range("B8").Value = -16.10616436612236
$ws.Range("F8").Value = 1.086010813713074
$ws.Range("J8").Value = -48.89927268028259
$ws.Range("K8").Value = -48.11897432804108
$ws.Range("B9").Value = -14.4238178851167
$ws.Range("F9").Value = 0.9476732015609741
$ws.Range("J9").Value = -49.38924181461334
$ws.Range("K9").Value = -49.02147305011749
$ws.Range("B10").Value = -11.46589212457798
$ws.Range("F10").Value = 0.1292674541473389
$ws.Range("J10").Value = -51.36069297790527
$ws.Range("K10").Value = -51.87237930297852
$ws.Range("B11").Value = -9.766161981559435
$ws.Range("F11").Value = -0.6806015968322754
$ws.Range("J11").Value = -53.98349905014038
$ws.Range("K11").Value = -54.9493613243103
$ws.Range("B12").Value = -9.375401323830147
$ws.Range("F12").Value = -0.9471328258514404
$ws.Range("J12").Value = -55.27644121646881
$ws.Range("K12").Value = -56.44125318527222
$ws.Range("B13").Value = -9.506691000926139
$ws.Range("F13").Value = -1.12610924243927
$ws.Range("J13").Value = -55.19175744056702
$ws.Range("K13").Value = -57.56403756141663
$ws.Range("B14").Value = -9.76961508686054
$ws.Range("F14").Value = -1.048430442810059
$ws.Range("J14").Value = -54.24193513393402
$ws.Range("K14").Value = -57.64688014984131
$ws.Range("B15").Value = -9.941163890847974
$ws.Range("F15").Value = -0.808897852897644
$ws.Range("J15").Value = -53.47002398967743
$ws.Range("K15").Value = -57.20019900798798
$ws.Range("B16").Value = -10.05301871215602
$ws.Range("F16").Value = -0.9695695638656616
$ws.Range("J16").Value = -53.33314919471741
$ws.Range("K16").Value = -57.3012284040451
$ws.Range("B17").Value = -10.16484558789693
$ws.Range("F17").Value = -1.183626055717468
$ws.Range("J17").Value = -53.26251769065857
$ws.Range("K17").Value = -57.4533154964447
$ws.Range("B18").Value = -10.46309461325495
$ws.Range("F18").Value = -0.9467613697052002
$ws.Range("J18").Value = -52.44827771186829
$ws.Range("K18").Value = -56.45097672939301
$ws.Range("B19").Value = -10.5806766469284
$ws.Range("F19").Value = -0.9132874011993408
$ws.Range("J19").Value = -52.28996860980988
$ws.Range("K19").Value = -55.9271422624588
$ws.Range("B20").Value = -9.942172894986015
$ws.Range("F20").Value = -0.9057048559188843
$ws.Range("J20").Value = -53.63475215435028
$ws.Range("K20").Value = -55.96534729003906
$ws.Range("B21").Value = -9.24502004647934
$ws.Range("F21").Value = -1.388172507286072
$ws.Range("J21").Value = -56.5741879940033
$ws.Range("K21").Value = -57.74955999851227
$ws.Range("B22").Value = -9.283099508225177
$ws.Range("F22").Value = -2.003084540367126
$ws.Range("J22").Value = -58.49073386192322
$ws.Range("K22").Value = -59.606818318367
$ws.Range("B23").Value = -9.488988295630179
$ws.Range("F23").Value = -2.561402678489685
$ws.Range("J23").Value = -59.97168600559235
$ws.Range("K23").Value = -61.12395787239075
$ws.Range("B24").Value = -9.600304425980539
$ws.Range("F24").Value = -2.858789205551147
$ws.Range("J24").Value = -60.54090583324432
$ws.Range("K24").Value = -61.75127637386322
$ws.Range("B25").Value = -9.34320657814942
$ws.Range("F25").Value = -2.422275543212891
$ws.Range("J25").Value = -57.7165333032608
$ws.Range("K25").Value = -59.21585619449615
